# Append new rows (171-185) of sweep-line filter-area results to the
# "result" worksheet, extending the used range from A1:D170 to A1:D185.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 28.40522003173828, 9.0625),
    @(99.67, 0.16, 34.94150567054749, 12.078125),
    @(100, 0.16, 37.01831841468811, 11.4375),
    @(100, 0.16, 34.95767211914062, 10.375),
    @(100, 0.16, 35.74127435684204, 10.84375),
    @(39.78, 0.16, 25.98540353775024, 9.796875),
    @(100, 0.16, 32.43182468414307, 10.359375),
    @(99.87, 0.16, 38.47451686859131, 14.4375),
    @(39.78, 0.16, 26.25908422470093, 9.65625),
    @(100, 0.16, 38.42415428161621, 12.078125),
    @(39.78, 0.16, 25.04115033149719, 8.1875),
    @(98.75, 0.16, 39.76158237457275, 13.71875),
    @(39.78, 0.16, 117.7494630813599, 10.5),
    @(39.78, 0.16, 80.7268750667572, 10.046875),
    @(100, 0.16, 37.93001079559326, 11.265625)
)

$startRow = 171
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
